$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 101
$ws.Range("A3").Value = 102
$ws.Range("A4").Value = 103
$ws.Range("A5").Value = 104

$ws.Range("A6").Select()
